$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gains a
#    trailing double-space, then three additional runs colored red are
#    appended: "(This is a change \u2013 Ve" / "rsion for main branch" / ")"
# ---------------------------------------------------------------------

$null = $d.Content.Find.Execute(
    "This is a Microsoft word document.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This is a Microsoft word document.  ", 2)

$p1 = $d.Paragraphs(1).Range
$insPos = $p1.End - 1
$ins1 = $d.Range($insPos, $insPos)
$ins1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$ins1.Font.Color = 255

$p1 = $d.Paragraphs(1).Range
$insPos = $p1.End - 1
$ins2 = $d.Range($insPos, $insPos)
$ins2.InsertAfter("rsion for main branch")
$ins2.Font.Color = 255

$p1 = $d.Paragraphs(1).Range
$insPos = $p1.End - 1
$ins3 = $d.Range($insPos, $insPos)
$ins3.InsertAfter(")")
$ins3.Font.Color = 255

# ---------------------------------------------------------------------
# 2) Drop the trailing "ank God almighty, we are free at last." paragraph
#    (the tail end of a word accidentally split across two paragraphs -
#    "Th" .. "ank God almighty..." - the whole paragraph is removed).
# ---------------------------------------------------------------------

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)
$lastPara.Range.Delete()

# ---------------------------------------------------------------------
# 3) Prune now-unused styles left over in styles.xml. Delete from the
#    highest style index down to the lowest to avoid shifting the
#    index of a not-yet-deleted style out from under a later lookup.
# ---------------------------------------------------------------------

$unusedStyleNames = @(
    "podcast-tools__subscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading 4 Char",
    "Heading 2 Char",
    "Hyperlink",
    "apple-converted-space",
    "Heading 4",
    "Heading 2"
)

foreach ($styleName in $unusedStyleNames) {
    $style = $d.Styles($styleName)
    $style.Delete()
}
